$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-04 01:39:39"

# Find last used row based on column A (header in row 1, data starts row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value -ne $null -and $cell.Value -ne "") {
        $cell.Value = $newTimestamp
    }
}
